# Revert csv module error handling — append the latest sensor reading
# (row 71) to each of the four per-row-lifter sheets, matching the
# pre-existing row layout (time, totals, hex payload, checksum, decoded
# numeric columns).

$wb = $excel.ActiveWorkbook

$newRow = 71

$data = @{
    "ROW35-FE-LIFTER"  = @{
        A = "2025-03-07 06:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    }
    "ROW35-MID-LIFTER" = @{
        A = "2025-03-07 06:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    }
    "ROW02-FE-LIFTER"  = @{
        A = "2025-03-07 06:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    }
    "ROW02-MID-LIFTER" = @{
        A = "2025-03-07 06:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $data[$sheetName]

    # Text columns (A-E): keep as plain text, same as every other data
    # row already on the sheet (inline/general text, no date or numeric
    # auto-conversion).
    $ws.Cells.Item($newRow, 1).Value = $row.A
    $ws.Cells.Item($newRow, 2).Value = $row.B
    $ws.Cells.Item($newRow, 3).Value = $row.C
    $ws.Cells.Item($newRow, 4).Value = $row.D
    $ws.Cells.Item($newRow, 5).Value = $row.E

    # Numeric columns (F, H, I)
    $ws.Cells.Item($newRow, 6).Value = $row.F
    $ws.Cells.Item($newRow, 8).Value = $row.H
    $ws.Cells.Item($newRow, 9).Value = $row.I

    # Column G is a long digit string that must stay text (Excel would
    # otherwise coerce it into a double and lose precision), so force a
    # text number format before assigning it.
    $gCell = $ws.Cells.Item($newRow, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row.G
}
